# Auto-generated cell value updates derived from the canonical OOXML diff.
# Each FFXIV Leve-profit worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) contains
# static market-price derived values (currentAveragePrice*, LevePrice*, LeveProfit*)
# in columns H-N that were refreshed by the scheduled data-update runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 468.97437
$ws.Range("J17").Value = 348.07895
$ws.Range("L17").Value = 1044.23685
$ws.Range("N17").Value = -1380.23685
$ws.Range("H43").Value = 6000.5
$ws.Range("I43").Value = 10001
$ws.Range("J43").Value = 2000
$ws.Range("K43").Value = 10001
$ws.Range("L43").Value = 2000
$ws.Range("M43").Value = -9932
$ws.Range("N43").Value = -2138
$ws.Range("H103").Value = 595
$ws.Range("I103").Value = 581.875
$ws.Range("J103").Value = 700
$ws.Range("K103").Value = 1745.625
$ws.Range("L103").Value = 2100
$ws.Range("M103").Value = -1159.625
$ws.Range("N103").Value = -3272
$ws.Range("H107").Value = 875
$ws.Range("I107").Value = 875
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 875
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1045
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 2986.4614
$ws.Range("I113").Value = 2383.8
$ws.Range("J113").Value = 3363.125
$ws.Range("K113").Value = 2383.8
$ws.Range("L113").Value = 3363.125
$ws.Range("M113").Value = 870.1999999999998
$ws.Range("N113").Value = -9871.125
$ws.Range("H127").Value = 922.1818
$ws.Range("J127").Value = 1956.25
$ws.Range("L127").Value = 5868.75
$ws.Range("N127").Value = -15788.75
$ws.Range("H129").Value = 1065.8903
$ws.Range("I129").Value = 340.6
$ws.Range("J129").Value = 1112.987
$ws.Range("K129").Value = 1021.8
$ws.Range("L129").Value = 3338.961
$ws.Range("M129").Value = 3978.2
$ws.Range("N129").Value = -13338.961
$ws.Range("H131").Value = 4531.607
$ws.Range("I131").Value = 1059.4445
$ws.Range("J131").Value = 10781.5
$ws.Range("K131").Value = 3178.3335
$ws.Range("L131").Value = 32344.5
$ws.Range("M131").Value = 1861.6665
$ws.Range("N131").Value = -42424.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9640.156000000001
$ws.Range("I32").Value = 9650.446
$ws.Range("J32").Value = 9555.556
$ws.Range("K32").Value = 9650.446
$ws.Range("L32").Value = 9555.556
$ws.Range("M32").Value = -9363.446
$ws.Range("N32").Value = -10129.556
$ws.Range("H88").Value = 2250.2856
$ws.Range("I88").Value = 1576
$ws.Range("K88").Value = 1576
$ws.Range("M88").Value = -1170
$ws.Range("H91").Value = 2250.2856
$ws.Range("I91").Value = 1576
$ws.Range("K91").Value = 1576
$ws.Range("M91").Value = -172
$ws.Range("H122").Value = 5390.8066
$ws.Range("I122").Value = 5784.6523
$ws.Range("J122").Value = 4258.5
$ws.Range("K122").Value = 17353.9569
$ws.Range("L122").Value = 12775.5
$ws.Range("M122").Value = -14903.9569
$ws.Range("N122").Value = -17675.5
$ws.Range("H135").Value = 26434.584
$ws.Range("J135").Value = 26434.584
$ws.Range("L135").Value = 26434.584
$ws.Range("N135").Value = -36574.584

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 87060.836
$ws.Range("I86").Value = 3784.8333
$ws.Range("J86").Value = 170336.83
$ws.Range("K86").Value = 3784.8333
$ws.Range("L86").Value = 170336.83
$ws.Range("M86").Value = -2661.8333
$ws.Range("N86").Value = -172582.83
$ws.Range("H89").Value = 87060.836
$ws.Range("I89").Value = 3784.8333
$ws.Range("J89").Value = 170336.83
$ws.Range("K89").Value = 18924.1665
$ws.Range("L89").Value = 851684.1499999999
$ws.Range("M89").Value = -13308.1665
$ws.Range("N89").Value = -862916.1499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1944.1167
$ws.Range("I31").Value = 2335.7036
$ws.Range("J31").Value = 1623.7273
$ws.Range("K31").Value = 2335.7036
$ws.Range("L31").Value = 1623.7273
$ws.Range("M31").Value = -2040.7036
$ws.Range("N31").Value = -2213.7273
$ws.Range("H34").Value = 1944.1167
$ws.Range("I34").Value = 2335.7036
$ws.Range("J34").Value = 1623.7273
$ws.Range("K34").Value = 2335.7036
$ws.Range("L34").Value = 1623.7273
$ws.Range("M34").Value = -2133.7036
$ws.Range("N34").Value = -2027.7273
$ws.Range("H99").Value = 1784
$ws.Range("I99").Value = 1772.4706
$ws.Range("J99").Value = 1980
$ws.Range("K99").Value = 1772.4706
$ws.Range("L99").Value = 1980
$ws.Range("M99").Value = -274.4706000000001
$ws.Range("N99").Value = -4976
$ws.Range("H107").Value = 754.7273
$ws.Range("I107").Value = 718.1875
$ws.Range("J107").Value = 852.1667
$ws.Range("K107").Value = 718.1875
$ws.Range("L107").Value = 852.1667
$ws.Range("M107").Value = 1201.8125
$ws.Range("N107").Value = -4692.1667
$ws.Range("H126").Value = 1784
$ws.Range("I126").Value = 1772.4706
$ws.Range("J126").Value = 1980
$ws.Range("K126").Value = 5317.4118
$ws.Range("L126").Value = 5940
$ws.Range("M126").Value = -2847.4118
$ws.Range("N126").Value = -10880

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 59500
$ws.Range("J37").Value = 59500
$ws.Range("L37").Value = 178500
$ws.Range("N37").Value = -178724
$ws.Range("H69").Value = 560
$ws.Range("I69").Value = 560
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 1680
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -869
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 560
$ws.Range("I72").Value = 560
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 5040
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -984
$ws.Range("N72").ClearContents()
$ws.Range("H75").Value = 5969.8667
$ws.Range("I75").Value = 1043.8
$ws.Range("J75").Value = 8432.9
$ws.Range("K75").Value = 3131.4
$ws.Range("L75").Value = 25298.7
$ws.Range("M75").Value = -2133.4
$ws.Range("N75").Value = -27294.7
$ws.Range("H78").Value = 5969.8667
$ws.Range("I78").Value = 1043.8
$ws.Range("J78").Value = 8432.9
$ws.Range("K78").Value = 9394.199999999999
$ws.Range("L78").Value = 75896.09999999999
$ws.Range("M78").Value = -4402.199999999999
$ws.Range("N78").Value = -85880.09999999999
$ws.Range("H107").Value = 1259.1428
$ws.Range("J107").Value = 1285.5172
$ws.Range("L107").Value = 3856.5516
$ws.Range("N107").Value = -7696.5516
$ws.Range("H121").Value = 37460.934
$ws.Range("I121").Value = 697.8889
$ws.Range("J121").Value = 92605.5
$ws.Range("K121").Value = 2093.6667
$ws.Range("L121").Value = 277816.5
$ws.Range("M121").Value = -783.6667000000002
$ws.Range("N121").Value = -280436.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2957.7812
$ws.Range("I102").Value = 2913.5
$ws.Range("J102").Value = 3149.6667
$ws.Range("K102").Value = 2913.5
$ws.Range("L102").Value = 3149.6667
$ws.Range("M102").Value = -1291.5
$ws.Range("N102").Value = -6393.6667
$ws.Range("H126").Value = 3185.889
$ws.Range("I126").Value = 1950
$ws.Range("J126").Value = 3539
$ws.Range("K126").Value = 5850
$ws.Range("L126").Value = 10617
$ws.Range("M126").Value = -3380
$ws.Range("N126").Value = -15557
$ws.Range("H141").Value = 56614.832
$ws.Range("J141").Value = 56614.832
$ws.Range("L141").Value = 56614.832
$ws.Range("N141").Value = -66974.83199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3700
$ws.Range("I40").Value = 3600
$ws.Range("J40").Value = 3900
$ws.Range("K40").Value = 3600
$ws.Range("L40").Value = 3900
$ws.Range("M40").Value = -3464
$ws.Range("N40").Value = -4172
$ws.Range("H46").Value = 1108.0714
$ws.Range("I46").Value = 939.375
$ws.Range("J46").Value = 1333
$ws.Range("K46").Value = 939.375
$ws.Range("L46").Value = 1333
$ws.Range("M46").Value = -751.375
$ws.Range("N46").Value = -1709
$ws.Range("H61").Value = 12001.579
$ws.Range("I61").Value = 14555.333
$ws.Range("J61").Value = 2425
$ws.Range("K61").Value = 14555.333
$ws.Range("L61").Value = 2425
$ws.Range("M61").Value = -14353.333
$ws.Range("N61").Value = -2829
$ws.Range("H93").Value = 1490.4
$ws.Range("I93").Value = 1266.6666
$ws.Range("K93").Value = 1266.6666
$ws.Range("M93").Value = -18.66660000000002
$ws.Range("H113").Value = 12001.579
$ws.Range("I113").Value = 14555.333
$ws.Range("J113").Value = 2425
$ws.Range("K113").Value = 14555.333
$ws.Range("L113").Value = 2425
$ws.Range("M113").Value = -12385.333
$ws.Range("N113").Value = -6765
$ws.Range("H132").Value = 3080.75
$ws.Range("I132").Value = 2933.1538
$ws.Range("K132").Value = 8799.4614
$ws.Range("M132").Value = -6269.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 82092.42999999999
$ws.Range("I81").Value = 113108.9
$ws.Range("J81").Value = 4551.25
$ws.Range("K81").Value = 226217.8
$ws.Range("L81").Value = 9102.5
$ws.Range("M81").Value = -225156.8
$ws.Range("N81").Value = -11224.5
$ws.Range("H84").Value = 82092.42999999999
$ws.Range("I84").Value = 113108.9
$ws.Range("J84").Value = 4551.25
$ws.Range("K84").Value = 1131089
$ws.Range("L84").Value = 45512.5
$ws.Range("M84").Value = -1125785
$ws.Range("N84").Value = -56120.5
$ws.Range("H100").Value = 38650
$ws.Range("I100").Value = 60720
$ws.Range("J100").Value = 1866.6666
$ws.Range("K100").Value = 121440
$ws.Range("L100").Value = 3733.3332
$ws.Range("M100").Value = -120899
$ws.Range("N100").Value = -4815.3332
$ws.Range("H122").Value = 73531064
$ws.Range("I122").Value = 104168170
$ws.Range("J122").Value = 2011.4
$ws.Range("K122").Value = 312504510
$ws.Range("L122").Value = 6034.200000000001
$ws.Range("M122").Value = -312502060
$ws.Range("N122").Value = -10934.2
$ws.Range("H126").Value = 9425.412
$ws.Range("I126").Value = 11787.385
$ws.Range("J126").Value = 1749
$ws.Range("K126").Value = 35362.155
$ws.Range("L126").Value = 5247
$ws.Range("M126").Value = -32892.155
$ws.Range("N126").Value = -10187
$ws.Range("H136").Value = 2313.8462
$ws.Range("I136").Value = 2203.0833
$ws.Range("J136").Value = 2491.0667
$ws.Range("K136").Value = 6609.249899999999
$ws.Range("L136").Value = 7473.2001
$ws.Range("M136").Value = -4059.249899999999
$ws.Range("N136").Value = -12573.2001
$ws.Range("H138").Value = 57137.95
$ws.Range("J138").Value = 57137.95
$ws.Range("L138").Value = 57137.95
$ws.Range("N138").Value = -67417.95
